# Apply changes described in the diff: add a new "2021" data column (R)
# to the worksheet, mirroring the formatting of the adjacent "2020" (Q)
# column, and move the active-cell selection to C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 0.12641839647678207
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 0.14922981985616976
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 0.10326895933792253
$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 0.03433011112114915
$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 0.036820478077087354
$ws.Range("Q9").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = 0.031930519190242035
$ws.Range("Q10").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = 0.087302929367211068
$ws.Range("Q11").Copy()
$ws.Range("R11").PasteSpecial(-4122)
$ws.Range("R11").Value = 0.10296328329317765
$ws.Range("Q12").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("R12").Value = 0.071859056271889668
$ws.Range("Q13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("R13").Value = 0.10716050460690947
$ws.Range("Q14").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$ws.Range("R14").Value = 0.079035451351703812
$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("R15").Value = 0.13553052227085377
$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("R16").Value = 0.06479643687803946
$ws.Range("Q17").Copy()
$ws.Range("R17").PasteSpecial(-4122)
$ws.Range("R17").Value = 0.07643825526207898
$ws.Range("Q18").Copy()
$ws.Range("R18").PasteSpecial(-4122)
$ws.Range("R18").Value = 0.053576570965516782
$ws.Range("Q19").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("R19").Value = 0.054163459619715498
$ws.Range("Q20").Copy()
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("R20").Value = 0.064872252119520635
$ws.Range("Q21").Copy()
$ws.Range("R21").PasteSpecial(-4122)
$ws.Range("R21").Value = 0.043693418784505472
$ws.Range("Q22").Copy()
$ws.Range("R22").PasteSpecial(-4122)
$ws.Range("R22").Value = 0.051373884452794741
$ws.Range("Q23").Copy()
$ws.Range("R23").PasteSpecial(-4122)
$ws.Range("R23").Value = 0.029662368095156877
$ws.Range("Q24").Copy()
$ws.Range("R24").PasteSpecial(-4122)
$ws.Range("R24").Value = 0.072642215296997686
$ws.Range("Q25").Copy()
$ws.Range("R25").PasteSpecial(-4122)
$ws.Range("R25").Value = 0.13772601093442507
$ws.Range("Q26").Copy()
$ws.Range("R26").PasteSpecial(-4122)
$ws.Range("R26").Value = 0.15668565643254884
$ws.Range("Q27").Copy()
$ws.Range("R27").PasteSpecial(-4122)
$ws.Range("R27").Value = 0.11816042869432726
$ws.Range("Q28").Copy()
$ws.Range("R28").PasteSpecial(-4122)
$ws.Range("R28").Value = 0.33417383115107696
$ws.Range("Q29").Copy()
$ws.Range("R29").PasteSpecial(-4122)
$ws.Range("R29").Value = 0.41139191068108794
$ws.Range("Q30").Copy()
$ws.Range("R30").PasteSpecial(-4122)
$ws.Range("R30").Value = 0.24697746624641295
$ws.Range("Q31").Copy()
$ws.Range("R31").PasteSpecial(-4122)
$ws.Range("R31").Value = 0.16773611144997194
$ws.Range("Q32").Copy()
$ws.Range("R32").PasteSpecial(-4122)
$ws.Range("R32").Value = 0.1959922553363346
$ws.Range("Q33").Copy()
$ws.Range("R33").PasteSpecial(-4122)
$ws.Range("R33").Value = 0.13791201213625709
$ws.Range("Q34").Copy()
$ws.Range("R34").PasteSpecial(-4122)
$ws.Range("Q35").Copy()
$ws.Range("R35").PasteSpecial(-4122)
$ws.Range("R35").Value = 0
$ws.Range("Q36").Copy()
$ws.Range("R36").PasteSpecial(-4122)
$ws.Range("R36").Value = 0.1
$ws.Range("Q37").Copy()
$ws.Range("R37").PasteSpecial(-4122)
$ws.Range("R37").Value = 0.2

# Update the selected / active cell shown when the workbook is reopened.
$ws.Range("C1").Select()
